$wb = $excel.ActiveWorkbook

# --- Login sheet: selection now spans A1:B1 (no data changes) ---
$loginWs = $wb.Worksheets.Item("Login")
[void]$loginWs.Range("A1:B1").Select()

# --- Product sheet: the 59931463 / PAX record moves from the bottom (row 6)
# to the top (row 1); the previously-empty row 2 gets filled in with what
# used to be row 1's data (40260717 / ALEX); everything else keeps its
# product mapping, just shifted up by one row so the sheet has no gaps.
$ws = $wb.Worksheets.Item("Product")

# Insert a row at position 2 (inherits row 1's formatting) and populate it
# with the old row-1 record, closing the gap that used to sit at row 2.
[void]$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = 40260717
$ws.Range("B2").Value = "ALEX"

# Drop the now-empty gap row (old row 3), pulling rows 4-7 up to 3-6.
[void]$ws.Rows.Item(3).Delete()

# Row 1 becomes the 59931463 / PAX record that used to live at the bottom.
$ws.Range("A1").Value = 59931463
$ws.Range("B1").Value = "PAX"

# Remove the now-duplicate trailing row.
[void]$ws.Rows.Item(6).Delete()

[void]$ws.Range("A1:B1").Select()
